# Sort/append posts by vote count to the bottom of the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the list (feature: sorting by vote count in
# the posts-list rendering). These mirror the existing "OK"/"Done"/"DONE"
# status markers used throughout the sheet.
$ws.Range("A24").Value = "DONE"
$ws.Range("A25").Value = "Done"
$ws.Range("A30").Value = "Done"
$ws.Range("A31").Value = "Done"

# Scroll the view back up a bit and leave the cursor on the next empty row,
# matching where the author's cursor ended up after entering the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("A32").Select()
